$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C9").Value = "'152"
$ws.Range("C9").Style = "Normal"
$ws.Range("D9").Value = "'399739.44"
$ws.Range("D9").Style = "Normal"

$ws.Range("C10").Value = "'324"
$ws.Range("C10").Style = "Normal"
$ws.Range("D10").Value = "'1150022.67"
$ws.Range("D10").Style = "Normal"

$ws.Range("C11").Value = "'138"
$ws.Range("C11").Style = "Normal"
$ws.Range("D11").Value = "'437451.40"
$ws.Range("D11").Style = "Normal"

$ws.Range("C14").Value = "'191"
$ws.Range("C14").Style = "Normal"
$ws.Range("D14").Value = "'477251.00"
$ws.Range("D14").Style = "Normal"

$ws.Range("C16").Value = "'450"
$ws.Range("C16").Style = "Normal"
$ws.Range("D16").Value = "'1558030.03"
$ws.Range("D16").Style = "Normal"

$ws.Range("C18").Value = "'14"
$ws.Range("C18").Style = "Normal"
$ws.Range("D18").Value = "'34971.00"
$ws.Range("D18").Style = "Normal"

$ws.Range("C20").Value = "'153"
$ws.Range("C20").Style = "Normal"
$ws.Range("D20").Value = "'390569.00"
$ws.Range("D20").Style = "Normal"

$ws.Range("C28").Value = "'209"
$ws.Range("C28").Style = "Normal"
$ws.Range("D28").Value = "'497573.00"
$ws.Range("D28").Style = "Normal"

$ws.Range("C30").Value = "'440"
$ws.Range("C30").Style = "Normal"
$ws.Range("D30").Value = "'1587507.82"
$ws.Range("D30").Style = "Normal"

$ws.Range("C32").Value = "'333"
$ws.Range("C32").Style = "Normal"
$ws.Range("D32").Value = "'1017882.96"
$ws.Range("D32").Style = "Normal"

$ws.Range("C44").Value = "'266"
$ws.Range("C44").Style = "Normal"
$ws.Range("D44").Value = "'732956.74"
$ws.Range("D44").Style = "Normal"

$ws.Range("C46").Value = "'515"
$ws.Range("C46").Style = "Normal"
$ws.Range("D46").Value = "'1805406.88"
$ws.Range("D46").Style = "Normal"

$ws.Range("C47").Value = "'334"
$ws.Range("C47").Style = "Normal"
$ws.Range("D47").Value = "'1065026.29"
$ws.Range("D47").Style = "Normal"

$ws.Range("C50").Value = "'2679"
$ws.Range("C50").Style = "Normal"
$ws.Range("D50").Value = "'5910081.17"
$ws.Range("D50").Style = "Normal"

$ws.Range("C51").Value = "'19"
$ws.Range("C51").Style = "Normal"
$ws.Range("D51").Value = "'108500.00"
$ws.Range("D51").Style = "Normal"

$ws.Range("C52").Value = "'3390"
$ws.Range("C52").Style = "Normal"
$ws.Range("D52").Value = "'10398583.25"
$ws.Range("D52").Style = "Normal"

$ws.Range("C54").Value = "'3507"
$ws.Range("C54").Style = "Normal"
$ws.Range("D54").Value = "'9984180.32"
$ws.Range("D54").Style = "Normal"

$ws.Range("C55").Value = "'50"
$ws.Range("C55").Style = "Normal"
$ws.Range("D55").Value = "'135350.00"
$ws.Range("D55").Style = "Normal"

$ws.Range("C56").Value = "'62"
$ws.Range("C56").Style = "Normal"
$ws.Range("D56").Value = "'188868.00"
$ws.Range("D56").Style = "Normal"

$ws.Range("C63").Value = "'208"
$ws.Range("C63").Style = "Normal"
$ws.Range("D63").Value = "'537013.00"
$ws.Range("D63").Style = "Normal"

$ws.Range("C64").Value = "'351"
$ws.Range("C64").Style = "Normal"
$ws.Range("D64").Value = "'1162870.80"
$ws.Range("D64").Style = "Normal"

$ws.Range("C65").Value = "'195"
$ws.Range("C65").Style = "Normal"
$ws.Range("D65").Value = "'597737.31"
$ws.Range("D65").Style = "Normal"

$ws.Range("C68").Value = "'345"
$ws.Range("C68").Style = "Normal"
$ws.Range("D68").Value = "'843135.70"
$ws.Range("D68").Style = "Normal"

$ws.Range("C70").Value = "'836"
$ws.Range("C70").Style = "Normal"
$ws.Range("D70").Value = "'2714081.34"
$ws.Range("D70").Style = "Normal"

$ws.Range("C71").Value = "'476"
$ws.Range("C71").Style = "Normal"
$ws.Range("D71").Value = "'1444224.03"
$ws.Range("D71").Style = "Normal"

$ws.Range("C73").Value = "'31"
$ws.Range("C73").Style = "Normal"
$ws.Range("D73").Value = "'113736.09"
$ws.Range("D73").Style = "Normal"

$ws.Range("C80").Value = "'192"
$ws.Range("C80").Style = "Normal"
$ws.Range("D80").Value = "'441071.00"
$ws.Range("D80").Style = "Normal"

$ws.Range("C82").Value = "'462"
$ws.Range("C82").Style = "Normal"
$ws.Range("D82").Value = "'1506362.50"
$ws.Range("D82").Style = "Normal"

$ws.Range("C83").Value = "'172"
$ws.Range("C83").Style = "Normal"
$ws.Range("D83").Value = "'483976.09"
$ws.Range("D83").Style = "Normal"

$ws.Range("C85").Value = "'6"
$ws.Range("C85").Style = "Normal"
$ws.Range("D85").Value = "'20670.00"
$ws.Range("D85").Style = "Normal"

$ws.Range("C86").Value = "'452"
$ws.Range("C86").Style = "Normal"
$ws.Range("D86").Value = "'1037224.67"
$ws.Range("D86").Style = "Normal"
